# Chỉnh sửa lại một vài lỗi
# Insert a new column before column A ("Mã lớp học phần") on Sheet1, shifting
# the existing Mã bài thi / Tên bài thi / ... table one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A; Excel shifts B:G (old A:F) to the right
# and carries along their formatting/column widths automatically.
$ws.Columns.Item(1).Insert()

# Fill in the new column's header + the three class-section codes.
$ws.Range("A1").Value = "Mã lớp học phần"
$ws.Range("A2").Value = "LHP_PTNC"
$ws.Range("A3").Value = "LHP_CSDL"
$ws.Range("A4").Value = "LHP_LTLN"

# Match the header/body cell styles used by the rest of the table
# (header style from B1, body style from B2) instead of the blank default.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# New column A gets its own (slightly narrower) width.
$ws.Columns.Item(1).ColumnWidth = 28.7

# Update the remembered selection to match the edited workbook.
$ws.Range("C9").Select() | Out-Null
